$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column before column N (14), shifting the
# "Late" / "Outstanding" columns one position to the right.
$ws.Columns.Item(14).Insert()

# New column N picks up column M's width (11.140625) but not its bestFit flag.
$ws.Columns.Item(14).ColumnWidth = $ws.Columns.Item(13).ColumnWidth

$ws.Select()
$ws.Range("R8").Select()
